# Apply the data-update edits described by the commit diff.
# The original workbook holds pause-frequency summary statistics; this
# revision refreshes a batch of mean/min-max/SD figures for the
# "disfluency" discipline rows (and a couple of neighbouring rows) and
# makes a small cosmetic tweak to the sheet view / column width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (A3=1, event-related / ls) -----------------------------------
$ws.Range("D3").Value = 52.18
$ws.Range("E3").Value = "10/139"
$ws.Range("F3").Value = 31.77
$ws.Range("G3").Value = 7.42
$ws.Range("H3").Value = "2.61/16.4"
$ws.Range("I3").Value = 2.9

# --- Row 7 (A7=5, disfluency / ls) ---------------------------------------
$ws.Range("D7").Value = 190.53
$ws.Range("E7").Value = "30/912"
$ws.Range("F7").Value = 166.54
$ws.Range("G7").Value = 26.65
$ws.Range("H7").Value = "5.06/105.37"
$ws.Range("I7").Value = 16.87

# --- Row 11 (A11=9, between-utterance / ls) ------------------------------
$ws.Range("D11").Value = 10.55
$ws.Range("F11").Value = 14.69
$ws.Range("G11").Value = 1.63
$ws.Range("I11").Value = 1.95

# --- Row 15 (A15=13, between-clause / ls) --------------------------------
$ws.Range("D15").Value = 402
$ws.Range("E15").Value = "60/1066"
$ws.Range("F15").Value = 205.75
$ws.Range("G15").Value = 56.92
$ws.Range("H15").Value = "27.19/108.19"
$ws.Range("I15").Value = 16.68

# --- Row 19 (A19=17, between-phrase / ls) --------------------------------
$ws.Range("D19").Value = 256.34
$ws.Range("E19").Value = "30/811"
$ws.Range("F19").Value = 156.88
$ws.Range("G19").Value = 36.24
$ws.Range("H19").Value = "13.59/78.85"
$ws.Range("I19").Value = 14.93

# --- Row 23 (A23=21, within-phrase / ls) ---------------------------------
$ws.Range("D23").Value = 120.26
$ws.Range("E23").Value = "13/335"
$ws.Range("F23").Value = 76.68
$ws.Range("G23").Value = 16.68
$ws.Range("H23").Value = "3.03/36.24"
$ws.Range("I23").Value = 7.25

# --- Cosmetic view tweaks ------------------------------------------------
# Column B was widened (best-fit) to accommodate the discipline labels.
$ws.Columns.Item(2).ColumnWidth = 14.6666666666667

# Scroll/selection state left by the author when the file was re-saved.
$ws.Activate() | Out-Null
$ws.Range("D23:I23").Select() | Out-Null
